# Auto-generated edit script: refreshes market-price-derived profit
# columns (H-N) on each class leve sheet to match a scheduled data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 86
$ws.Range("H86").Value = 3264.3333
$ws.Range("I86").Value = 1980.4
$ws.Range("J86").Value = 4181.4287
$ws.Range("K86").Value = 1980.4
$ws.Range("L86").Value = 4181.4287
$ws.Range("M86").Value = -857.4000000000001
$ws.Range("N86").Value = -6427.4287
# Row 89
$ws.Range("H89").Value = 3264.3333
$ws.Range("I89").Value = 1980.4
$ws.Range("J89").Value = 4181.4287
$ws.Range("K89").Value = 9902
$ws.Range("L89").Value = 20907.1435
$ws.Range("M89").Value = -4286
$ws.Range("N89").Value = -32139.1435
# Row 116
$ws.Range("H116").Value = 41178.332
$ws.Range("I116").Value = 52811.2
$ws.Range("J116").Value = 21790.223
$ws.Range("K116").Value = 52811.2
$ws.Range("L116").Value = 21790.223
$ws.Range("M116").Value = -49369.2
$ws.Range("N116").Value = -28674.223
# Row 138
$ws.Range("H138").Value = 3137.35
$ws.Range("I138").Value = 2207.9565
$ws.Range("J138").Value = 3715.081
$ws.Range("K138").Value = 6623.869499999999
$ws.Range("L138").Value = 11145.243
$ws.Range("M138").Value = -1483.869499999999
$ws.Range("N138").Value = -21425.243
# Row 141
$ws.Range("I141").Value = 1096.1111
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 3288.3333
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 1891.6667
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 25
$ws.Range("H25").Value = 1471.25
$ws.Range("J25").Value = 2500
$ws.Range("L25").Value = 2500
$ws.Range("N25").Value = -3304
# Row 32
$ws.Range("H32").Value = 15219.833
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 15219.833
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 15219.833
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -15793.833
# Row 45
$ws.Range("H45").Value = 4280.8667
$ws.Range("I45").Value = 3434.4167
$ws.Range("J45").Value = 7666.6665
$ws.Range("K45").Value = 3434.4167
$ws.Range("L45").Value = 7666.6665
$ws.Range("M45").Value = -3057.4167
$ws.Range("N45").Value = -8420.666499999999
# Row 74
$ws.Range("H74").Value = 3714.1936
$ws.Range("I74").Value = 3140.347
$ws.Range("K74").Value = 3140.347
$ws.Range("M74").Value = -2266.347
# Row 77
$ws.Range("H77").Value = 3714.1936
$ws.Range("I77").Value = 3140.347
$ws.Range("K77").Value = 15701.735
$ws.Range("M77").Value = -11333.735
# Row 102
$ws.Range("H102").Value = 4247.636
$ws.Range("I102").Value = 4035.4119
$ws.Range("J102").Value = 4969.2
$ws.Range("K102").Value = 4035.4119
$ws.Range("L102").Value = 4969.2
$ws.Range("M102").Value = -2413.4119
$ws.Range("N102").Value = -8213.200000000001
# Row 110
$ws.Range("H110").Value = 2723.55
$ws.Range("I110").Value = 2198.4443
$ws.Range("K110").Value = 2198.4443
$ws.Range("M110").Value = -153.4443000000001
# Row 132
$ws.Range("H132").Value = 4085.8604
$ws.Range("I132").Value = 2546.1965
$ws.Range("K132").Value = 7638.5895
$ws.Range("M132").Value = -5108.5895
# Row 135
$ws.Range("H135").Value = 62887.11
$ws.Range("J135").Value = 62887.11
$ws.Range("L135").Value = 62887.11
$ws.Range("N135").Value = -73027.11

$ws = $wb.Worksheets.Item("BSM")
# Row 80
$ws.Range("H80").Value = 5973.923
$ws.Range("J80").Value = 3070.625
$ws.Range("L80").Value = 3070.625
$ws.Range("N80").Value = -5066.625
# Row 83
$ws.Range("H83").Value = 5973.923
$ws.Range("J83").Value = 3070.625
$ws.Range("L83").Value = 15353.125
$ws.Range("N83").Value = -25337.125
# Row 99
$ws.Range("H99").Value = 3943.8823
$ws.Range("I99").Value = 2343.5
$ws.Range("J99").Value = 5366.4443
$ws.Range("K99").Value = 2343.5
$ws.Range("L99").Value = 5366.4443
$ws.Range("M99").Value = -845.5
$ws.Range("N99").Value = -8362.444299999999
# Row 105
$ws.Range("H105").Value = 2996.4333
$ws.Range("I105").Value = 2662.524
$ws.Range("J105").Value = 3775.5557
$ws.Range("K105").Value = 2662.524
$ws.Range("L105").Value = 3775.5557
$ws.Range("M105").Value = -915.5239999999999
$ws.Range("N105").Value = -7269.5557
# Row 107
$ws.Range("H107").Value = 2189.9443
$ws.Range("I107").Value = 1877.2413
$ws.Range("J107").Value = 3485.4285
$ws.Range("K107").Value = 1877.2413
$ws.Range("L107").Value = 3485.4285
$ws.Range("M107").Value = 42.75870000000009
$ws.Range("N107").Value = -7325.4285

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 9822.946
$ws.Range("I31").Value = 2821.2
$ws.Range("K31").Value = 2821.2
$ws.Range("M31").Value = -2526.2
# Row 34
$ws.Range("H34").Value = 9822.946
$ws.Range("I34").Value = 2821.2
$ws.Range("K34").Value = 2821.2
$ws.Range("M34").Value = -2619.2
# Row 58
$ws.Range("H58").Value = 4251
$ws.Range("I58").Value = 3603
$ws.Range("K58").Value = 3603
$ws.Range("M58").Value = -3400
# Row 99
$ws.Range("H99").Value = 4293.2144
$ws.Range("I99").Value = 2482.8333
$ws.Range("J99").Value = 5651
$ws.Range("K99").Value = 2482.8333
$ws.Range("L99").Value = 5651
$ws.Range("M99").Value = -984.8332999999998
$ws.Range("N99").Value = -8647
# Row 105
$ws.Range("H105").Value = 1734
$ws.Range("I105").Value = 1734
$ws.Range("K105").Value = 1734
$ws.Range("M105").Value = 13
# Row 126
$ws.Range("H126").Value = 4293.2144
$ws.Range("I126").Value = 2482.8333
$ws.Range("J126").Value = 5651
$ws.Range("K126").Value = 7448.499899999999
$ws.Range("L126").Value = 16953
$ws.Range("M126").Value = -4978.499899999999
$ws.Range("N126").Value = -21893
# Row 132
$ws.Range("H132").Value = 3287.2856
$ws.Range("I132").Value = 1714.7916
$ws.Range("K132").Value = 5144.3748
$ws.Range("M132").Value = -2614.3748
# Row 134
$ws.Range("H134").Value = 9818.853999999999
$ws.Range("I134").Value = 9818.853999999999
$ws.Range("K134").Value = 29456.562
$ws.Range("M134").Value = -26921.562
# Row 136
$ws.Range("H136").Value = 4251
$ws.Range("I136").Value = 3603
$ws.Range("K136").Value = 10809
$ws.Range("M136").Value = -8259

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 166688.47
$ws.Range("J2").Value = 42.25
$ws.Range("L2").Value = 253.5
$ws.Range("N2").Value = -479.5
# Row 6
$ws.Range("H6").Value = 153
$ws.Range("I6").Value = 136.66667
$ws.Range("J6").Value = 300
$ws.Range("K6").Value = 410.00001
$ws.Range("L6").Value = 900
$ws.Range("M6").Value = -297.00001
$ws.Range("N6").Value = -1126
# Row 86
$ws.Range("H86").Value = 1096.7778
$ws.Range("J86").Value = 1134.2
$ws.Range("L86").Value = 3402.6
$ws.Range("N86").Value = -5774.6
# Row 89
$ws.Range("H89").Value = 1096.7778
$ws.Range("J89").Value = 1134.2
$ws.Range("L89").Value = 10207.8
$ws.Range("N89").Value = -22063.8

$ws = $wb.Worksheets.Item("GSM")
# Row 62
$ws.Range("H62").Value = 75063.75
$ws.Range("J62").Value = 75063.75
$ws.Range("L62").Value = 75063.75
$ws.Range("N62").Value = -76435.75
# Row 65
$ws.Range("H65").Value = 75063.75
$ws.Range("J65").Value = 75063.75
$ws.Range("L65").Value = 225191.25
$ws.Range("N65").Value = -232055.25
# Row 70
$ws.Range("H70").Value = 5748.643
$ws.Range("I70").Value = 5474
$ws.Range("J70").Value = 5901.222
$ws.Range("K70").Value = 5474
$ws.Range("L70").Value = 5901.222
$ws.Range("M70").Value = -5204
$ws.Range("N70").Value = -6441.222
# Row 73
$ws.Range("H73").Value = 5748.643
$ws.Range("I73").Value = 5474
$ws.Range("J73").Value = 5901.222
$ws.Range("K73").Value = 5474
$ws.Range("L73").Value = 5901.222
$ws.Range("M73").Value = -4538
$ws.Range("N73").Value = -7773.222
# Row 80
$ws.Range("H80").Value = 165126.03
$ws.Range("I80").Value = 266700.75
$ws.Range("J80").Value = 4299.4165
$ws.Range("K80").Value = 266700.75
$ws.Range("L80").Value = 4299.4165
$ws.Range("M80").Value = -265702.75
$ws.Range("N80").Value = -6295.4165
# Row 83
$ws.Range("H83").Value = 165126.03
$ws.Range("I83").Value = 266700.75
$ws.Range("J83").Value = 4299.4165
$ws.Range("K83").Value = 1333503.75
$ws.Range("L83").Value = 21497.0825
$ws.Range("M83").Value = -1328511.75
$ws.Range("N83").Value = -31481.0825
# Row 97
$ws.Range("H97").Value = 1201.75
$ws.Range("I97").Value = 533.9375
$ws.Range("J97").Value = 1869.5625
$ws.Range("K97").Value = 533.9375
$ws.Range("L97").Value = 1869.5625
$ws.Range("M97").Value = -37.9375
$ws.Range("N97").Value = -2861.5625
# Row 99
$ws.Range("H99").Value = 5946.091
$ws.Range("I99").Value = 901.1429000000001
$ws.Range("K99").Value = 901.1429000000001
$ws.Range("M99").Value = 1344.8571
# Row 126
$ws.Range("H126").Value = 3015.087
$ws.Range("I126").Value = 2713
$ws.Range("K126").Value = 8139
$ws.Range("M126").Value = -5669
# Row 132
$ws.Range("H132").Value = 2732.1785
$ws.Range("I132").Value = 2181.2727
$ws.Range("K132").Value = 6543.8181
$ws.Range("M132").Value = -4013.8181

$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 3994.0454
$ws.Range("J68").Value = 4451.6665
$ws.Range("L68").Value = 4451.6665
$ws.Range("N68").Value = -5949.6665
# Row 71
$ws.Range("H71").Value = 3994.0454
$ws.Range("J71").Value = 4451.6665
$ws.Range("L71").Value = 22258.3325
$ws.Range("N71").Value = -29746.3325
# Row 93
$ws.Range("H93").Value = 2261.2
$ws.Range("I93").Value = 1764.0625
$ws.Range("J93").Value = 4249.75
$ws.Range("K93").Value = 1764.0625
$ws.Range("L93").Value = 4249.75
$ws.Range("M93").Value = -516.0625
$ws.Range("N93").Value = -6745.75
# Row 122
$ws.Range("H122").Value = 3530.9092
$ws.Range("I122").Value = 3327.625
$ws.Range("K122").Value = 9982.875
$ws.Range("M122").Value = -7532.875
# Row 132
$ws.Range("H132").Value = 3834.1924
$ws.Range("I132").Value = 3134.1
$ws.Range("K132").Value = 9402.299999999999
$ws.Range("M132").Value = -6872.299999999999

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 5400.8
$ws.Range("I122").Value = 5157.25
$ws.Range("K122").Value = 15471.75
$ws.Range("M122").Value = -13021.75
# Row 132
$ws.Range("H132").Value = 11880837
$ws.Range("I132").Value = 1284.375
$ws.Range("K132").Value = 3853.125
$ws.Range("M132").Value = -1323.125
